$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
Write-Output (Get-Member -InputObject $s | Select-String "Theme" | Out-String)
